$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-detected as a number by Excel;
# force them to remain Text, matching the workbook's existing text cells.
$textCells = @("D5","D8","D11","D14","D15","D18","D19","D22","D25","D26","D28","D30","D32","D41","D42","D43","D49","D50")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "29.847.52"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.626.55"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.72%  "
$ws.Range("D5").Value = "214.38"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").Value = "28.46"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").Value = "0.0899"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "1.860.47"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "1.626.27"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "0.563"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "9.22"
$ws.Range("E15").Value = "  +4.40%  "
$ws.Range("D16").Value = "29.856.48"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "64.55"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "239.82"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "0.0₃0701"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "9.78"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").Value = "157.63"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").Value = "15.44"
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "6.54"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").Value = "0.0488"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").Value = "1.422.57"
$ws.Range("E35").Value = "  +3.79%  "
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").Value = "74.61"
$ws.Range("E41").Value = "  +7.15%  "
$ws.Range("D42").Value = "0.0499"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").Value = "0.827"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "1.766.86"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("D49").Value = "48.35"
$ws.Range("E49").Value = "  -10.99%  "
$ws.Range("D50").Value = "90.88"
$ws.Range("E50").Value = "  +3.83%  "
$ws.Range("E51").Value = "  +8.57%  "

# Restore default (Normal) style on the forced-text cells so no stray style index remains.
foreach ($c in $textCells) { $ws.Range($c).Style = "Normal" }
